$wb = $excel.ActiveWorkbook

# Update DatosCuenta (sheet1) values for the "Catorce" smoke row -> "Uno"
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "SmokeUno"
$wsCuenta.Range("B2").Value = "SmokeLastNUno"
$wsCuenta.Range("B3").Select()

# Make DatosMotor the active/selected sheet
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Select()
